# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" sheet (fund holdings) before the "总计" sheet,
#    formatted like the existing quarterly sheets (2021-Q1..Q3).
# 2. Insert a new top row into "总计" summarizing the 2022-Q1 quarter and
#    bump the existing row index counters down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: add the "2022-Q1" worksheet right after "2021-Q3" (i.e. right
# before "总计"), using "2021-Q3" as the formatting template.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q3")
$newSheet = $wb.Worksheets.Add($null, $template)
$newSheet.Name = "2022-Q1"

# Clone cell formatting from the template (header + 4 data rows) so the
# new sheet matches the look of the other quarterly sheets.
$template.Range("B1:H4").Copy()
$newSheet.Range("B1").PasteSpecial(-4122)
$wb.Application.CutCopyMode = $false

$template.Range("A2:A4").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)
$wb.Application.CutCopyMode = $false

$template.Range("A4").Copy()
$newSheet.Range("A5").PasteSpecial(-4122)
$wb.Application.CutCopyMode = $false

$template.Range("B4:H4").Copy()
$newSheet.Range("B5").PasteSpecial(-4122)
$wb.Application.CutCopyMode = $false

# Header row
$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

# Fund holdings data for 2022-Q1
$funds = @(
  @{code="561550"; name="华泰柏瑞中证500增强策略ETF";   scale="9.68"; pos="98.93"; pct="1.09"; mv="0.1055"; rank=8},
  @{code="008114"; name="天弘中证红利低波动100指数A";     scale="3.16"; pos="92.60"; pct="1.59"; mv="0.0502"; rank=10},
  @{code="008115"; name="天弘中证红利低波动100指数C";     scale="2.37"; pos="92.60"; pct="1.59"; mv="0.0377"; rank=10},
  @{code="515100"; name="景顺长城中证红利低波动100ETF";   scale="1.25"; pos="97.96"; pct="1.68"; mv="0.0210"; rank=10}
)

$r = 2
foreach ($fund in $funds) {
  # Columns that hold numeric-looking text must be forced to Text format
  # so leading/trailing zeros in values like "008114" / "92.60" survive.
  $newSheet.Range("B$r").NumberFormat = "@"
  $newSheet.Range("D$r`:G$r").NumberFormat = "@"

  $newSheet.Cells.Item($r, 1).Value = $r - 2
  $newSheet.Cells.Item($r, 2).Value = $fund.code
  $newSheet.Cells.Item($r, 3).Value = $fund.name
  $newSheet.Cells.Item($r, 4).Value = $fund.scale
  $newSheet.Cells.Item($r, 5).Value = $fund.pos
  $newSheet.Cells.Item($r, 6).Value = $fund.pct
  $newSheet.Cells.Item($r, 7).Value = $fund.mv
  $newSheet.Cells.Item($r, 8).Value = $fund.rank
  $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 2: update the "总计" (summary) sheet with a new first data row for
# 2022-Q1, shifting the previous rows down by one.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()

# Carry the index-column style down onto the newly inserted A2 cell.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$wb.Application.CutCopyMode = $false

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q1"
$summary.Cells.Item(2, 3).Value = 4
$summary.Cells.Item(2, 4).Value = 0.21

# Renumber the index column for the rows that shifted down.
$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(5, 1).Value = 3
